$wb = $excel.ActiveWorkbook

# --- Update text/date values across sheets ---
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# "2016-09-07 17:26:49" -> "2016-09-07 17:27:48" (Overview G2, de-de H2)
$overview.Range("G2").Value = "2016-09-07 17:27:48"
$dede.Range("H2").Value = "2016-09-07 17:27:48"

# "2016-09-07 17:26:44" -> "2016-09-07 17:27:42" (zh-cn H2)
$zhcn.Range("H2").Value = "2016-09-07 17:27:42"

# --- Update column widths ---
$overview.Range("E1:F1").ColumnWidth = 17.2159881591797
$zhcn.Range("C1").ColumnWidth = 17.2159881591797
$dede.Range("C1").ColumnWidth = 17.2159881591797
